$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Columns("D:E").Insert()

# Copy number formats/fonts from column F (old column D) into new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restated financial data for all quarters (columns D through M)
$data = @{
    7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    8 = @(581300, 749800, 1002400, 780000, 979100, 636900, 629700, 574000, 854700, 675800)
    9 = @(224400, 391600, 573900, 512200, 631200, 321800, 321700, 304900, 466700, 360900)
    10 = @(356900, 358200, 428500, 267700, 347900, 315200, 308000, 269000, 388000, 315000)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(481200, 562900, 748400, 654000, 816600, 466500, 479700, 475500, 619100, 496200)
    18 = @(100000, 187000, 254000, 126000, 162400, 170400, 150000, 98400, 235600, 179600)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(33800, 10300, 4700, 18300, 40700, 40700, 31500, 13600, 7200, 16500)
    21 = @(178000, 237400, 298900, 181200, 242300, 251700, 221600, 150700, 290500, 234000)
    22 = @(30400, 18400, 18800, 20800, 43100, 40200, 42100, 20300, "NA", 21100)
    23 = @(103400, 178800, 239900, 123400, 160000, 171000, 139500, 91700, 242800, 175100)
    24 = @(27400, 39700, 66600, 24700, 14900, 28700, 23100, 11400, 46700, 26100)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(76100, 139100, 173300, 98700, 145100, 142300, 116400, 80400, 196100, 149000)
    27 = @(72100, 119300, 151000, 59000, 112800, 108100, 79400, 70000, 173800, 124700)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(-33800, -10300, -4700, -18300, -40700, -40700, -31500, -13600, -7200, -16500)
    33 = @(72100, 119300, 151000, 59000, 112800, 108100, 79400, 70000, 173800, 124700)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(72100, 119300, 151000, 59000, 112800, 108100, 79400, 70000, 173800, 124700)
    38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(1688000, 1928800, 1999600, 2512700, 2784200, 2582400, 2429900, 2567300, 2690100, 2103600)
    42 = @(10500, 10800, 11600, 11400, 11600, 12100, 13000, 13000, 12000, 11700)
    43 = @(783600, 790000, 804500, 1025600, 888100, 711700, 896000, 851100, 854300, 1019200)
    44 = @(4215600, 4332300, 4293500, 3218600, 3184500, 3730100, 3830700, 3812300, 3823600, 4508500)
    45 = @(11700, 2800, 3700, 2800, 54000, 220900, 2800, 2800, 2900, 202000)
    46 = @(6709300, 7064600, 7113000, 6771100, 6922300, 7257200, 7172400, 7246600, 7382900, 7845000)
    47 = @(1931800, 1765700, 1764700, 1730400, 1426600, 1393500, 1370100, 1368200, 1362300, 1175600)
    48 = @(6455200, 5766600, 5519300, 5516400, 5491600, 5573600, 5511600, 5411500, 5479600, 5451000)
    49 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(303700, 455400, 550000, 435100, 437700, 423100, 271400, 271000, 274600, 112600)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(15400100, 15052300, 14947000, 14453000, 14278300, 14647400, 14325500, 14297200, 14499400, 14584200)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(953600, 1015500, 1077800, 1059100, 958000, 1424300, 1202900, 1124800, 1153700, 1301500)
    58 = @(927900, 853800, 727200, 743100, 933500, 1274000, 912100, 1063400, 1305700, 1361200)
    59 = @(407000, 283500, 325000, 301900, 554500, 249700, 195000, 204500, 220700, 193400)
    60 = @(2288600, 2152800, 2130000, 2104200, 2446100, 2948000, 2310000, 2392700, 2680100, 2856100)
    61 = @(3737500, 3307300, 3258000, 2813500, 2769200, 2552800, 3033500, 3021400, 2896600, 3059300)
    62 = @(323800, 443100, 425200, 459100, 475700, 526900, 565200, 552700, 567200, 505200)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(7996600, 7555400, 7484100, 7030400, 7353600, 7685900, 7452300, 7483300, 7692700, 7983600)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 243500, 0, 0, 0, 241800, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(6038700, 6101700, 6032500, 5979800, 5513500, 5891100, 5797600, 5763200, 5689000, 5562200)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(7403500, 7496900, 7462900, 7422500, 6681200, 6961500, 6873200, 6813900, 6564800, 6600600)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    81 = @(72100, 119300, 151000, 59000, 112800, 108100, 79400, 70000, 173800, 124700)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(44100, 40200, 40100, 36900, 39300, 40700, 40100, 38900, 47700, 37900)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(204700, 148600, -768300, -27200, 370200, 122000, 223700, 85400, 302800, 389000)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(-51000, -32700, -25000, -39600, -22100, -14300, -19100, -49900, -36300, -35700)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-934600, -287300, -28600, -24000, -48000, 80000, -82100, -12800, 647800, -330600)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(-7700, -66700, -110600, -25300, -4400, -48300, -99900, -27100, -8100, -46800)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(477300, 4800, 400100, -220400, -122300, -31400, -289900, -234600, -468700, -305800)
    101 = @(2900, -15300, 600, 7500, -4900, -5600, 5200, -17500, 15600, 4000)
    102 = @(-249800, -149200, -396300, -264100, 202800, 171400, -144400, -204600, 497500, -243500)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item([int]$row, 4 + $i).Value = $vals[$i]
        }
    }
}

Write-Host "Edit complete"